$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": add column CA (31-aug) ---
$ws1 = $wb.Worksheets.Item("Prix Spot")

$ws1.Range("CA1").Value = "31-aug"
$ws1.Range("CA2").Value = 15.86
$ws1.Range("CA3").Value = 12.93
$ws1.Range("CA4").Value = 12.33
$ws1.Range("CA5").Value = 7.61
$ws1.Range("CA6").Value = 5.95
$ws1.Range("CA7").Value = 8.12
$ws1.Range("CA8").Value = 9.18
$ws1.Range("CA9").Value = 16.24
$ws1.Range("CA10").Value = 22.15
$ws1.Range("CA11").Value = 3.6
$ws1.Range("CA12").Value = 0
$ws1.Range("CA13").Value = -0.01
$ws1.Range("CA14").Value = -0.08
$ws1.Range("CA15").Value = -0.95
$ws1.Range("CA16").Value = -0.6
$ws1.Range("CA17").Value = -0.01
$ws1.Range("CA18").Value = 4.22
$ws1.Range("CA19").Value = 4.96
$ws1.Range("CA20").Value = 17.53
$ws1.Range("CA21").Value = 47.15
$ws1.Range("CA22").Value = 59.55
$ws1.Range("CA23").Value = 73.88
$ws1.Range("CA24").Value = 71.4
$ws1.Range("CA25").Value = 62

# Match header style of neighbouring header cell (BZ1): bold, centered, bordered
$ws1.Range("BZ1").Copy()
$ws1.Range("CA1").PasteSpecial(-4122)

# --- Sheet "Gaz": add row 76 ---
$ws2 = $wb.Worksheets.Item("Gaz")
# Force text type so the date-like string isn't auto-converted to a date value,
# then restore the plain (unstyled) look of the existing data rows.
$ws2.Range("A76").NumberFormat = "@"
$ws2.Range("A76").Value = "2025-08-29"
$ws2.Range("A75").Copy()
$ws2.Range("A76").PasteSpecial(-4122)
$ws2.Range("B76").Value = 30.375

# --- Sheet "CO2": add row 76 ---
$ws3 = $wb.Worksheets.Item("CO2")
$ws3.Range("A76").NumberFormat = "@"
$ws3.Range("A76").Value = "2025-08-29"
$ws3.Range("A75").Copy()
$ws3.Range("A76").PasteSpecial(-4122)
$ws3.Range("B76").Value = 71.1
